# Updates cryptos list values per the Wed Sep  6 16:45:32 UTC 2023 GitHub Actions run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.701.28"
$ws.Range("E2").Value = "  -0.93%  "
$ws.Range("D3").Value = "1.630.11"
$ws.Range("E3").Value = "  -1.05%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.04"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.09%  "
$ws.Range("E6").Value = "  -1.07%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  -0.88%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0636"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.42%  "
$ws.Range("E10").Value = "  -5.77%  "
$ws.Range("E11").Value = "  -0.07%  "
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.24"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.08%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.631.14"
$ws.Range("E13").Value = "  -0.93%  "
$ws.Range("D14").Value = "1.854.44"
$ws.Range("E14").Value = "  -1.02%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.552"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.02%  "
$ws.Range("D16").Value = "0.0₃0764"
$ws.Range("E16").Value = "  -1.20%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.09"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.64%  "
$ws.Range("D18").Value = "25.719.09"
$ws.Range("E18").Value = "  -0.90%  "
$ws.Range("E20").Value = "  +0.80%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "193.61"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.33%  "
$ws.Range("E22").Value = "  -0.50%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.19"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.84%  "
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.78"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.22%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "140.02"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.04%  "
$ws.Range("E27").Value = "  -3.24%  "
$ws.Range("E28").Value = "  -0.51%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.48"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.66%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.24"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.91%  "
$ws.Range("E31").Value = "  -2.77%  "
$ws.Range("E32").Value = "  +0.39%  "
$ws.Range("E33").Value = "  -0.09%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.59"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.30%  "
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("E36").Value = "  -1.64%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.55"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.67%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.543"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.69%  "
$ws.Range("D39").Value = "1.103.21"
$ws.Range("E39").Value = "  -2.64%  "
$ws.Range("E40").Value = "  -1.33%  "
$ws.Range("E41").Value = "  -0.16%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "99.74"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.95%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.794"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.15%  "
$ws.Range("D45").Value = "1.762.64"
$ws.Range("E45").Value = "  -1.07%  "
$ws.Range("D46").Value = "0.0₆0109"
$ws.Range("E46").Value = "  -2.00%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "54.98"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.65%  "
$ws.Range("E48").Value = "  -2.60%  "
$ws.Range("B49").Value = "SynthetixNetwork"
$ws.Range("C49").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.37"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.34%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0503"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.42%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.65"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.64%  "
